$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-7, columns M:T
$data = @{
    2 = @{ M=45.1830845;          N=90.366169;          O=0.2982772948921854; P=0.2359735829156887; Q=1.264237765338167;  R=7.585426592029;    S=0.2982772948921854; T=0.2359735829156887 }
    3 = @{ M=24.489942;           N=73.46982600000001;  O=0.1616709822417395; P=0.1918520865636367; Q=0.6852367404740001; R=6.167130664266002; S=0.1616709822417395; T=0.1918520865636367 }
    4 = @{ M=20.755341;           N=62.26602299999999;  O=0.1370169176485697; P=0.1625955454769879; Q=0.5807413596269999; R=5.226672236642999; S=0.1370169176485697; T=0.1625955454769879 }
    5 = @{ M=26.3069545;          N=52.613909;          O=0.173666037012409;  P=0.1373909368441856; Q=0.7360773558948334; R=4.416464135369;    S=0.173666037012409;  T=0.1373909368441856 }
    6 = @{ M=17.34473466666667;   N=52.034204;          O=0.1145017121838161; P=0.1358771505744131; Q=0.4853114575515556; R=4.367803117964001; S=0.1145017121838161; T=0.1358771505744131 }
    7 = @{ M=17.400077;           N=52.200231;          O=0.1148670560212801; P=0.136310697625088;  Q=0.4868599544856667; R=4.381739590371;    S=0.1148670560212801; T=0.136310697625088  }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("M$row").Value = $cols.M
    $ws.Range("N$row").Value = $cols.N
    $ws.Range("O$row").Value = $cols.O
    $ws.Range("P$row").Value = $cols.P
    $ws.Range("Q$row").Value = $cols.Q
    $ws.Range("R$row").Value = $cols.R
    $ws.Range("S$row").Value = $cols.S
    $ws.Range("T$row").Value = $cols.T
}
